$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.975.01'
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").Value = '1.643.01'
$ws.Range("E3").Value = '  -0.45%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.70%  '

$ws.Range("D5").Value = '215.31'
$ws.Range("E5").Value = '  -0.14%  '

$ws.Range("D6").Value = '0.5059'
$ws.Range("E6").Value = '  -1.04%  '

$ws.Range("E7").Value = '  -0.72%  '

$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("D9").Value = '0.06417'
$ws.Range("E9").Value = '  -0.32%  '

$ws.Range("D10").Value = '19.71'
$ws.Range("E10").Value = '  +0.00%  '

$ws.Range("D11").Value = '0.07757'
$ws.Range("E11").Value = '  +0.35%  '

$ws.Range("D12").Value = '4.276'
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("D13").Value = '1.648.60'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").Value = '1.869.93'
$ws.Range("E14").Value = '  -0.45%  '

$ws.Range("D15").Value = '0.5471'
$ws.Range("E15").Value = '  -0.20%  '

$ws.Range("D16").Value = '0.0₅7938'
$ws.Range("E16").Value = '  -0.80%  '

$ws.Range("D17").Value = '64.46'
$ws.Range("E17").Value = '  +0.91%  '

$ws.Range("D18").Value = '25.981.41'
$ws.Range("E18").Value = '  -0.54%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.73%  '

$ws.Range("D20").Value = '203.07'
$ws.Range("E20").Value = '  -2.14%  '

$ws.Range("D21").Value = '4.395'
$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("D22").Value = '9.908'
$ws.Range("E22").Value = '  -1.34%  '

$ws.Range("D23").Value = '5.993'

$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.75%  '

$ws.Range("D25").Value = '1.875'
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '140.86'
$ws.Range("E26").Value = '  -1.55%  '

$ws.Range("D27").Value = '0.1138'
$ws.Range("E27").Value = '  -2.65%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '15.73'
$ws.Range("E28").Value = '  -0.54%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '6.825'
$ws.Range("E29").Value = '  -1.36%  '

$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("D31").Value = '0.04938'
$ws.Range("E31").Value = '  -2.75%  '

$ws.Range("E32").Value = '  -2.20%  '

$ws.Range("D33").Value = '3.211'
$ws.Range("E33").Value = '  -0.69%  '

$ws.Range("E34").Value = '  -0.31%  '

$ws.Range("D35").Value = '2.361'
$ws.Range("E35").Value = '  +0.50%  '

$ws.Range("D36").Value = '0.8943'
$ws.Range("E36").Value = '  -2.33%  '

$ws.Range("D37").Value = '2.621'
$ws.Range("E37").Value = '  -0.66%  '

$ws.Range("D38").Value = '1.150.36'
$ws.Range("E38").Value = '  +0.53%  '

$ws.Range("D39").Value = '0.5589'
$ws.Range("E39").Value = '  -1.92%  '

$ws.Range("E40").Value = '  -0.45%  '

$ws.Range("E41").Value = '  -0.78%  '

$ws.Range("D42").Value = '5.708'

$ws.Range("D43").Value = '0.8074'
$ws.Range("E43").Value = '  -2.15%  '

$ws.Range("D44").Value = '99.74'
$ws.Range("E44").Value = '  -0.37%  '

$ws.Range("D45").Value = '1.780.73'

$ws.Range("D46").Value = '0.0₈118'
$ws.Range("E46").Value = '  +4.82%  '

$ws.Range("D47").Value = '0.4515'
$ws.Range("E47").Value = '  -0.51%  '

$ws.Range("D48").Value = '1.003'
$ws.Range("E48").Value = '  -0.83%  '

$ws.Range("D49").Value = '54.79'
$ws.Range("E49").Value = '  -0.80%  '

$ws.Range("D50").Value = '0.05044'
$ws.Range("E50").Value = '  -0.61%  '

$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  -0.69%  '
